$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C8").Value = 111
$wb.Save()
